$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data collected for row-groups 41-91 (M/N/O/P columns): raw counts,
# average (N), a newly-collected O metric, and P = percentage of total ---

# Row 41
$ws.Range("M41").Value = 22
$ws.Range("N41").Formula = "=AVERAGE(B41:M41)"
$ws.Range("O41").Value = 15
$ws.Range("P41").Formula = "=O41*100/31"
# Row 42
$ws.Range("M42").Value = 2
$ws.Range("N42").Formula = "=AVERAGE(B42:M42)"
$ws.Range("O42").Value = 7
$ws.Range("P42").Formula = "=O42*100/31"
# Row 43
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = 0
$ws.Range("N43").Formula = "=AVERAGE(B43:M43)"
$ws.Range("O43").Value = 0
$ws.Range("P43").Formula = "=O43*100/31"
# Row 44
$ws.Range("C44").Formula = "=31-C43-C42-C41"
$ws.Range("D44").Formula = "=31-D43-D42-D41"
$ws.Range("E44").Formula = "=31-E43-E42-E41"
$ws.Range("G44").Formula = "=31-G43-G42-G41"
$ws.Range("H44").Formula = "=31-H43-H42-H41"
$ws.Range("I44").Formula = "=31-I43-I42-I41"
$ws.Range("J44").Formula = "=31-J43-J42-J41"
$ws.Range("K44").Formula = "=31-K43-K42-K41"
$ws.Range("L44").Formula = "=31-L43-L42-L41"
$ws.Range("M44").Formula = "=31-M43-M42-M41"
$ws.Range("N44").Formula = "=31-N43-N42-N41"
$ws.Range("O44").Formula = "=31-O43-O42-O41"
$ws.Range("P44").Formula = "=O44*100/31"
# Row 49
$ws.Range("M49").Value = 7
$ws.Range("N49").Formula = "=AVERAGE(B49:M49)"
$ws.Range("O49").Value = 6
$ws.Range("P49").Formula = "=O49*100/26"
# Row 50
$ws.Range("M50").Value = 3
$ws.Range("N50").Formula = "=AVERAGE(B50:M50)"
$ws.Range("O50").Value = 5
$ws.Range("P50").Formula = "=O50*100/26"
# Row 51
$ws.Range("E51").Value = 0
$ws.Range("H51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("M51").Value = 0
$ws.Range("N51").Formula = "=AVERAGE(B51:M51)"
$ws.Range("O51").Value = 1
$ws.Range("P51").Formula = "=O51*100/26"
# Row 52
$ws.Range("C52").Formula = "=26-C51-C50-C49"
$ws.Range("D52").Formula = "=26-D51-D50-D49"
$ws.Range("E52").Formula = "=26-E51-E50-E49"
$ws.Range("G52").Formula = "=26-G51-G50-G49"
$ws.Range("H52").Formula = "=26-H51-H50-H49"
$ws.Range("I52").Formula = "=26-I51-I50-I49"
$ws.Range("J52").Formula = "=26-J51-J50-J49"
$ws.Range("K52").Formula = "=26-K51-K50-K49"
$ws.Range("L52").Formula = "=26-L51-L50-L49"
$ws.Range("M52").Formula = "=26-M51-M50-M49"
$ws.Range("N52").Formula = "=26-N51-N50-N49"
$ws.Range("O52").Formula = "=26-O51-O50-O49"
$ws.Range("P52").Formula = "=O52*100/26"
# Row 57
$ws.Range("M57").Value = 4
$ws.Range("N57").Formula = "=AVERAGE(B57:M57)"
$ws.Range("O57").Value = 4
$ws.Range("P57").Formula = "=O57*100/17"
# Row 58
$ws.Range("M58").Value = 1
$ws.Range("N58").Formula = "=AVERAGE(B58:M58)"
$ws.Range("O58").Value = 2
$ws.Range("P58").Formula = "=O58*100/17"
# Row 59
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = 2
$ws.Range("N59").Formula = "=AVERAGE(B59:M59)"
$ws.Range("O59").Value = 1
$ws.Range("P59").Formula = "=O59*100/17"
# Row 60
$ws.Range("C60").Formula = "=17-C59-C58-C57"
$ws.Range("D60").Formula = "=17-D59-D58-D57"
$ws.Range("E60").Formula = "=17-E59-E58-E57"
$ws.Range("G60").Formula = "=17-G59-G58-G57"
$ws.Range("H60").Formula = "=17-H59-H58-H57"
$ws.Range("I60").Formula = "=17-I59-I58-I57"
$ws.Range("J60").Formula = "=17-J59-J58-J57"
$ws.Range("K60").Formula = "=17-K59-K58-K57"
$ws.Range("L60").Formula = "=17-L59-L58-L57"
$ws.Range("M60").Formula = "=17-M59-M58-M57"
$ws.Range("N60").Formula = "=17-N59-N58-N57"
$ws.Range("O60").Formula = "=17-O59-O58-O57"
$ws.Range("P60").Formula = "=O60*100/17"
# Row 64
$ws.Range("M64").Value = 9
$ws.Range("N64").Formula = "=AVERAGE(B64:M64)"
$ws.Range("O64").Value = 6
$ws.Range("P64").Formula = "=O64*100/24"
# Row 65
$ws.Range("M65").Value = 6
$ws.Range("N65").Formula = "=AVERAGE(B65:M65)"
$ws.Range("O65").Value = 6
$ws.Range("P65").Formula = "=O65*100/24"
# Row 66
$ws.Range("E66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = 1
$ws.Range("N66").Formula = "=AVERAGE(B66:M66)"
$ws.Range("O66").Value = 0
$ws.Range("P66").Formula = "=O66*100/24"
# Row 67
$ws.Range("C67").Formula = "=24-C66-C65-C64"
$ws.Range("D67").Formula = "=24-D66-D65-D64"
$ws.Range("E67").Formula = "=24-E66-E65-E64"
$ws.Range("G67").Formula = "=24-G66-G65-G64"
$ws.Range("H67").Formula = "=24-H66-H65-H64"
$ws.Range("I67").Formula = "=24-I66-I65-I64"
$ws.Range("J67").Formula = "=24-J66-J65-J64"
$ws.Range("K67").Formula = "=24-K66-K65-K64"
$ws.Range("L67").Formula = "=24-L66-L65-L64"
$ws.Range("M67").Formula = "=24-M66-M65-M64"
$ws.Range("N67").Formula = "=24-N66-N65-N64"
$ws.Range("O67").Formula = "=24-O66-O65-O64"
$ws.Range("P67").Formula = "=O67*100/24"
# Row 72
$ws.Range("M72").Value = 4
$ws.Range("N72").Formula = "=AVERAGE(B72:M72)"
$ws.Range("O72").Value = 3
$ws.Range("P72").Formula = "=O72*100/11"
# Row 73
$ws.Range("M73").Value = 2
$ws.Range("N73").Formula = "=AVERAGE(B73:M73)"
$ws.Range("O73").Value = 3
$ws.Range("P73").Formula = "=O73*100/11"
# Row 74
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 1
$ws.Range("N74").Formula = "=AVERAGE(B74:M74)"
$ws.Range("O74").Value = 0
$ws.Range("P74").Formula = "=O74*100/11"
# Row 75
$ws.Range("C75").Formula = "=11-C74-C73-C72"
$ws.Range("D75").Formula = "=11-D74-D73-D72"
$ws.Range("E75").Formula = "=11-E74-E73-E72"
$ws.Range("F75").Formula = "=11-F74-F73-F72"
$ws.Range("G75").Formula = "=11-G74-G73-G72"
$ws.Range("H75").Formula = "=11-H74-H73-H72"
$ws.Range("I75").Formula = "=11-I74-I73-I72"
$ws.Range("J75").Formula = "=11-J74-J73-J72"
$ws.Range("K75").Formula = "=11-K74-K73-K72"
$ws.Range("L75").Formula = "=11-L74-L73-L72"
$ws.Range("M75").Formula = "=11-M74-M73-M72"
$ws.Range("N75").Formula = "=11-N74-N73-N72"
$ws.Range("O75").Formula = "=11-O74-O73-O72"
$ws.Range("P75").Formula = "=O75*100/11"
# Row 80
$ws.Range("M80").Value = 9
$ws.Range("N80").Formula = "=AVERAGE(B80:M80)"
$ws.Range("O80").Value = 7
$ws.Range("P80").Formula = "=O80*100/19"
# Row 81
$ws.Range("M81").Value = 3
$ws.Range("N81").Formula = "=AVERAGE(B81:M81)"
$ws.Range("O81").Value = 3
$ws.Range("P81").Formula = "=O81*100/19"
# Row 82
$ws.Range("D82").Value = 0
$ws.Range("E82").Value = 0
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").Value = 0
$ws.Range("N82").Formula = "=AVERAGE(B82:M82)"
$ws.Range("O82").Value = 0
$ws.Range("P82").Formula = "=O82*100/19"
# Row 83
$ws.Range("C83").Formula = "=19-C82-C81-C80"
$ws.Range("D83").Formula = "=19-D82-D81-D80"
$ws.Range("E83").Formula = "=19-E82-E81-E80"
$ws.Range("F83").Formula = "=19-F82-F81-F80"
$ws.Range("G83").Formula = "=19-G82-G81-G80"
$ws.Range("H83").Formula = "=19-H82-H81-H80"
$ws.Range("I83").Formula = "=19-I82-I81-I80"
$ws.Range("J83").Formula = "=19-J82-J81-J80"
$ws.Range("K83").Formula = "=19-K82-K81-K80"
$ws.Range("L83").Formula = "=19-L82-L81-L80"
$ws.Range("M83").Formula = "=19-M82-M81-M80"
$ws.Range("N83").Formula = "=19-N82-N81-N80"
$ws.Range("O83").Formula = "=19-O82-O81-O80"
$ws.Range("P83").Formula = "=O83*100/19"
# Row 88
$ws.Range("M88").Value = 8
$ws.Range("N88").Formula = "=AVERAGE(B88:M88)"
$ws.Range("O88").Value = 8
$ws.Range("P88").Formula = "=O88*100/67"
# Row 89
$ws.Range("M89").Value = 16
$ws.Range("N89").Formula = "=AVERAGE(B89:M89)"
$ws.Range("O89").Value = 12
$ws.Range("P89").Formula = "=O89*100/67"
# Row 90
$ws.Range("F90").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = 0
$ws.Range("N90").Formula = "=AVERAGE(B90:M90)"
$ws.Range("O90").Value = 1
$ws.Range("P90").Formula = "=O90*100/67"
# Row 91
$ws.Range("C91").Formula = "=67-C90-C89-C88"
$ws.Range("D91").Formula = "=67-D90-D89-D88"
$ws.Range("E91").Formula = "=67-E90-E89-E88"
$ws.Range("F91").Formula = "=67-F90-F89-F88"
$ws.Range("G91").Formula = "=67-G90-G89-G88"
$ws.Range("H91").Formula = "=67-H90-H89-H88"
$ws.Range("I91").Formula = "=67-I90-I89-I88"
$ws.Range("J91").Formula = "=67-J90-J89-J88"
$ws.Range("K91").Formula = "=67-K90-K89-K88"
$ws.Range("L91").Formula = "=67-L90-L89-L88"
$ws.Range("M91").Formula = "=67-M90-M89-M88"
$ws.Range("N91").Formula = "=67-N90-N89-N88"
$ws.Range("O91").Formula = "=67-O90-O89-O88"
$ws.Range("P91").Formula = "=O91*100/67"

# --- View-state: update the selected range / active cell to match the
# latest edit location ---
$ws.Range("P88:P91").Select()
